# Auto-generated Excel COM-interop script applying the cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.587.57'
$ws.Range('E2').Value = '  +3.09%  '

$ws.Range('D3').Value = '1.904.53'
$ws.Range('E3').Value = '  +1.07%  '

$ws.Range('D4').Value = '''0.9999'
$ws.Range('D4').ClearFormats()
$ws.Range('D4').PrefixCharacter = ''
$ws.Range('E4').Value = '  -0.93%  '

$ws.Range('D5').Value = '''314.75'
$ws.Range('D5').ClearFormats()
$ws.Range('D5').PrefixCharacter = ''
$ws.Range('E5').Value = '  -0.36%  '

$ws.Range('D6').Value = '''1.000'
$ws.Range('D6').ClearFormats()
$ws.Range('D6').PrefixCharacter = ''
$ws.Range('E6').Value = '  -0.72%  '

$ws.Range('D7').Value = '''0.5153'
$ws.Range('D7').ClearFormats()
$ws.Range('D7').PrefixCharacter = ''
$ws.Range('E7').Value = '  +0.59%  '

$ws.Range('D8').Value = '''0.3947'
$ws.Range('D8').ClearFormats()
$ws.Range('D8').PrefixCharacter = ''
$ws.Range('E8').Value = '  -0.41%  '

$ws.Range('D9').Value = '''0.08453'
$ws.Range('D9').ClearFormats()
$ws.Range('D9').PrefixCharacter = ''
$ws.Range('E9').Value = '  +0.12%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '''1.120'
$ws.Range('D10').ClearFormats()
$ws.Range('D10').PrefixCharacter = ''
$ws.Range('E10').Value = '  +0.31%  '

$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '''42.39'
$ws.Range('D11').ClearFormats()
$ws.Range('D11').PrefixCharacter = ''
$ws.Range('E11').Value = '  +1.48%  '

$ws.Range('D12').Value = '''6.291'
$ws.Range('D12').ClearFormats()
$ws.Range('D12').PrefixCharacter = ''
$ws.Range('E12').Value = '  +0.03%  '

$ws.Range('D13').Value = '1.904.32'
$ws.Range('E13').Value = '  +1.33%  '

$ws.Range('D14').Value = '''20.76'
$ws.Range('D14').ClearFormats()
$ws.Range('D14').PrefixCharacter = ''
$ws.Range('E14').Value = '  +0.88%  '

$ws.Range('D15').Value = '''7.328'
$ws.Range('D15').ClearFormats()
$ws.Range('D15').PrefixCharacter = ''
$ws.Range('E15').Value = '  +0.56%  '

$ws.Range('D16').Value = '''0.9998'
$ws.Range('D16').ClearFormats()
$ws.Range('D16').PrefixCharacter = ''
$ws.Range('E16').Value = '  -1.07%  '

$ws.Range('D17').Value = '''93.09'
$ws.Range('D17').ClearFormats()
$ws.Range('D17').PrefixCharacter = ''
$ws.Range('E17').Value = '  +1.76%  '

$ws.Range('D18').Value = '''0.00001109'
$ws.Range('D18').ClearFormats()
$ws.Range('D18').PrefixCharacter = ''
$ws.Range('E18').Value = '  -0.06%  '

$ws.Range('D19').Value = '''0.06752'
$ws.Range('D19').ClearFormats()
$ws.Range('D19').PrefixCharacter = ''
$ws.Range('E19').Value = '  -0.06%  '

$ws.Range('D20').Value = '''17.92'
$ws.Range('D20').ClearFormats()
$ws.Range('D20').PrefixCharacter = ''
$ws.Range('E20').Value = '  +0.82%  '

$ws.Range('D21').Value = '''1.002'
$ws.Range('D21').ClearFormats()
$ws.Range('D21').PrefixCharacter = ''
$ws.Range('E21').Value = '  -0.51%  '

$ws.Range('D22').Value = '''6.019'
$ws.Range('D22').ClearFormats()
$ws.Range('D22').PrefixCharacter = ''
$ws.Range('E22').Value = '  +0.72%  '

$ws.Range('D23').Value = '29.554.61'
$ws.Range('E23').Value = '  +2.93%  '

$ws.Range('D24').Value = '''11.18'
$ws.Range('D24').ClearFormats()
$ws.Range('D24').PrefixCharacter = ''
$ws.Range('E24').Value = '  +0.30%  '

$ws.Range('D25').Value = '''2.212'
$ws.Range('D25').ClearFormats()
$ws.Range('D25').PrefixCharacter = ''
$ws.Range('E25').Value = '  -1.83%  '

$ws.Range('D26').Value = '2.112.65'
$ws.Range('E26').Value = '  +0.88%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''20.98'
$ws.Range('D27').ClearFormats()
$ws.Range('D27').PrefixCharacter = ''
$ws.Range('E27').Value = '  +0.54%  '

$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''158.85'
$ws.Range('D28').ClearFormats()
$ws.Range('D28').PrefixCharacter = ''
$ws.Range('E28').Value = '  -1.53%  '

$ws.Range('D29').Value = '''2.441'
$ws.Range('D29').ClearFormats()
$ws.Range('D29').PrefixCharacter = ''
$ws.Range('E29').Value = '  +2.06%  '

$ws.Range('D30').Value = '''127.89'
$ws.Range('D30').ClearFormats()
$ws.Range('D30').PrefixCharacter = ''
$ws.Range('E30').Value = '  +0.18%  '

$ws.Range('D31').Value = '''1.065'
$ws.Range('D31').ClearFormats()
$ws.Range('D31').PrefixCharacter = ''
$ws.Range('E31').Value = '  +1.02%  '

$ws.Range('D32').Value = '''0.1050'
$ws.Range('D32').ClearFormats()
$ws.Range('D32').PrefixCharacter = ''
$ws.Range('E32').Value = '  -0.52%  '

$ws.Range('D33').Value = '''6.160'
$ws.Range('D33').ClearFormats()
$ws.Range('D33').PrefixCharacter = ''
$ws.Range('E33').Value = '  +5.98%  '

$ws.Range('D34').Value = '''3.650'
$ws.Range('D34').ClearFormats()
$ws.Range('D34').PrefixCharacter = ''
$ws.Range('E34').Value = '  +0.74%  '

$ws.Range('D35').Value = '''0.02490'
$ws.Range('D35').ClearFormats()
$ws.Range('D35').PrefixCharacter = ''
$ws.Range('E35').Value = '  +1.06%  '

$ws.Range('D36').Value = '''0.06602'
$ws.Range('D36').ClearFormats()
$ws.Range('D36').PrefixCharacter = ''
$ws.Range('E36').Value = '  +0.88%  '

$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').Value = '''9.073'
$ws.Range('D37').ClearFormats()
$ws.Range('D37').PrefixCharacter = ''
$ws.Range('E37').Value = '  +1.24%  '

$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '''0.2202'
$ws.Range('D38').ClearFormats()
$ws.Range('D38').PrefixCharacter = ''
$ws.Range('E38').Value = '  +0.29%  '

$ws.Range('D39').Value = '''5.231'
$ws.Range('D39').ClearFormats()
$ws.Range('D39').PrefixCharacter = ''
$ws.Range('E39').Value = '  +2.55%  '

$ws.Range('D40').Value = '''1.235'
$ws.Range('D40').ClearFormats()
$ws.Range('D40').PrefixCharacter = ''
$ws.Range('E40').Value = '  +2.52%  '

$ws.Range('D41').Value = '''0.6545'
$ws.Range('D41').ClearFormats()
$ws.Range('D41').PrefixCharacter = ''
$ws.Range('E41').Value = '  +1.18%  '

$ws.Range('E42').Value = '  -2.30%  '

$ws.Range('D43').Value = '''11.30'
$ws.Range('D43').ClearFormats()
$ws.Range('D43').PrefixCharacter = ''
$ws.Range('E43').Value = '  +0.71%  '

$ws.Range('D44').Value = '''0.6086'
$ws.Range('D44').ClearFormats()
$ws.Range('D44').PrefixCharacter = ''
$ws.Range('E44').Value = '  +0.01%  '

$ws.Range('D45').Value = '''13.23'
$ws.Range('D45').ClearFormats()
$ws.Range('D45').PrefixCharacter = ''

$ws.Range('E46').Value = '  -0.87%  '

$ws.Range('D47').Value = '''2.059'
$ws.Range('D47').ClearFormats()
$ws.Range('D47').PrefixCharacter = ''
$ws.Range('E47').Value = '  +1.91%  '

$ws.Range('D48').Value = '''1.233'
$ws.Range('D48').ClearFormats()
$ws.Range('D48').PrefixCharacter = ''
$ws.Range('E48').Value = '  +1.35%  '

$ws.Range('D49').Value = '''123.83'
$ws.Range('D49').ClearFormats()
$ws.Range('D49').PrefixCharacter = ''
$ws.Range('E49').Value = '  +0.99%  '

$ws.Range('D50').Value = '''1.160'
$ws.Range('D50').ClearFormats()
$ws.Range('D50').PrefixCharacter = ''
$ws.Range('E50').Value = '  -2.82%  '

$ws.Range('D51').Value = '''77.95'
$ws.Range('D51').ClearFormats()
$ws.Range('D51').PrefixCharacter = ''
$ws.Range('E51').Value = '  +0.74%  '
